$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rtn4"
$ws.Range("C2").Value = "Tnfrsf19"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 40.797777
$ws.Range("H2").Value = 122.393331
$ws.Range("I2").Value = 0.2689231481273683
$ws.Range("J2").Value = 0.2689231481273683
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05206533333333333
$ws.Range("N2").Value = 0.156196
$ws.Range("O2").Value = 0.03170654174267026
$ws.Range("P2").Value = 0.03170654174267026
$ws.Range("Q2").Value = 2.124149858764
$ws.Range("R2").Value = 19.117348728876
$ws.Range("S2").Value = 0.008526623021670699
$ws.Range("T2").Value = 0.008526623021670702

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rtn4"
$ws.Range("C3").Value = "Tnfrsf19"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 40.797777
$ws.Range("H3").Value = 122.393331
$ws.Range("I3").Value = 0.2689231481273683
$ws.Range("J3").Value = 0.2689231481273683
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.150436
$ws.Range("N3").Value = 0.451308
$ws.Range("O3").Value = 0.09161192310175054
$ws.Range("P3").Value = 0.09161192310175056
$ws.Range("Q3").Value = 6.137454380772
$ws.Range("R3").Value = 55.237089426948
$ws.Range("S3").Value = 0.02463656676652513
$ws.Range("T3").Value = 0.02463656676652514

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rtn4"
$ws.Range("C4").Value = "Tnfrsf19"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 40.797777
$ws.Range("H4").Value = 122.393331
$ws.Range("I4").Value = 0.2689231481273683
$ws.Range("J4").Value = 0.2689231481273683
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.439599333333333
$ws.Range("N4").Value = 4.318798
$ws.Range("O4").Value = 0.8766815351555791
$ws.Range("P4").Value = 0.8766815351555792
$ws.Range("Q4").Value = 58.73245257068201
$ws.Range("R4").Value = 528.5920731361381
$ws.Range("S4").Value = 0.2357599583391724
$ws.Range("T4").Value = 0.2357599583391725

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rtn4"
$ws.Range("C5").Value = "Tnfrsf19"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 46.219831
$ws.Range("H5").Value = 138.659493
$ws.Range("I5").Value = 0.3046632285488233
$ws.Range("J5").Value = 0.3046632285488233
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05206533333333333
$ws.Range("N5").Value = 0.156196
$ws.Range("O5").Value = 0.03170654174267026
$ws.Range("P5").Value = 0.03170654174267026
$ws.Range("Q5").Value = 2.406450907625333
$ws.Range("R5").Value = 21.658058168628
$ws.Range("S5").Value = 0.009659817373439955
$ws.Range("T5").Value = 0.009659817373439957

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rtn4"
$ws.Range("C6").Value = "Tnfrsf19"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 46.219831
$ws.Range("H6").Value = 138.659493
$ws.Range("I6").Value = 0.3046632285488233
$ws.Range("J6").Value = 0.3046632285488233
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.150436
$ws.Range("N6").Value = 0.451308
$ws.Range("O6").Value = 0.09161192310175054
$ws.Range("P6").Value = 0.09161192310175056
$ws.Range("Q6").Value = 6.953126496315999
$ws.Range("R6").Value = 62.578138466844
$ws.Range("S6").Value = 0.02791078426574585
$ws.Range("T6").Value = 0.02791078426574586

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rtn4"
$ws.Range("C7").Value = "Tnfrsf19"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 46.219831
$ws.Range("H7").Value = 138.659493
$ws.Range("I7").Value = 0.3046632285488233
$ws.Range("J7").Value = 0.3046632285488233
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.439599333333333
$ws.Range("N7").Value = 4.318798
$ws.Range("O7").Value = 0.8766815351555791
$ws.Range("P7").Value = 0.8766815351555792
$ws.Range("Q7").Value = 66.53803789437934
$ws.Range("R7").Value = 598.842341049414
$ws.Range("S7").Value = 0.2670926269096375
$ws.Range("T7").Value = 0.2670926269096375

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Rtn4"
$ws.Range("C8").Value = "Tnfrsf19"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 64.69033266666666
$ws.Range("H8").Value = 194.070998
$ws.Range("I8").Value = 0.4264136233238083
$ws.Range("J8").Value = 0.4264136233238083
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.05206533333333333
$ws.Range("N8").Value = 0.156196
$ws.Range("O8").Value = 0.03170654174267026
$ws.Range("P8").Value = 0.03170654174267026
$ws.Range("Q8").Value = 3.368123733734222
$ws.Range("R8").Value = 30.313113603608
$ws.Range("S8").Value = 0.0135201013475596
$ws.Range("T8").Value = 0.0135201013475596

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Rtn4"
$ws.Range("C9").Value = "Tnfrsf19"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 64.69033266666666
$ws.Range("H9").Value = 194.070998
$ws.Range("I9").Value = 0.4264136233238083
$ws.Range("J9").Value = 0.4264136233238083
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.150436
$ws.Range("N9").Value = 0.451308
$ws.Range("O9").Value = 0.09161192310175054
$ws.Range("P9").Value = 0.09161192310175056
$ws.Range("Q9").Value = 9.731754885042665
$ws.Range("R9").Value = 87.585793965384
$ws.Range("S9").Value = 0.03906457206947955
$ws.Range("T9").Value = 0.03906457206947956

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Rtn4"
$ws.Range("C10").Value = "Tnfrsf19"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 64.69033266666666
$ws.Range("H10").Value = 194.070998
$ws.Range("I10").Value = 0.4264136233238083
$ws.Range("J10").Value = 0.4264136233238083
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.439599333333333
$ws.Range("N10").Value = 4.318798
$ws.Range("O10").Value = 0.8766815351555791
$ws.Range("P10").Value = 0.8766815351555792
$ws.Range("Q10").Value = 93.1281597800449
$ws.Range("R10").Value = 838.153438020404
$ws.Range("S10").Value = 0.3738289499067691
$ws.Range("T10").Value = 0.3738289499067692
